$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 - header values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 - Subject "CON" (B2:E2)
$ws.Range("B2").Value = 34.584652283769003
$ws.Range("C2").Value = 19.784121002394567
$ws.Range("D2").Value = 44.367993887109812
$ws.Range("E2").Value = 18.762756017546945

# Row 3 - Subject "STR" (B3:E3)
$ws.Range("B3").Value = 37.532249289626257
$ws.Range("C3").Value = 19.421284614683866
$ws.Range("D3").Value = 54.57503551642754
$ws.Range("E3").Value = 24.705470356675846

# Update selection to match the new active selection range B1:E3
$ws.Range("B1:E3").Select()
